# "Actualizar" automation run: refresh the last-checked timestamp on the
# most recent availability cycle (rows 366-379) and append a brand new
# 14-row availability-check cycle (rows 380-393) with the same
# Nombre/URL/Disponibilidad pattern used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Touch up the timestamp of the previous cycle (rows 366-379): the
#    automation re-saved the same check with a marginally later instant.
# ---------------------------------------------------------------------
$refreshedTimestamp = 44232.21908939815
for ($r = 366; $r -le 379; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $refreshedTimestamp
}

# ---------------------------------------------------------------------
# 2) Append the new cycle: rows 380-393.
# ---------------------------------------------------------------------
$names = @("Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE", "UtilidadesOdoo", "Filtros Dashboard", "MapStore", "GeoServer", "Tomcat", "Shiny", "Github", "EZ Exporter")
# Hyperlink target (rels "Target" - no fragment).
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# Only the MapStore link (index 8) carries a "/" fragment sub-address,
# matching every earlier MapStore row in the sheet.
$subAddresses = @("", "", "", "", "", "", "", "", "/", "", "", "", "", "")
# Column B's visible text is the full URL the cell shows, which for the
# MapStore row includes the "#/" fragment (matches every earlier row).
$cellTexts = @()
for ($i = 0; $i -lt $urls.Count; $i++) {
    if ($subAddresses[$i] -ne "") {
        $cellTexts += ($urls[$i] + "#" + $subAddresses[$i])
    } else {
        $cellTexts += $urls[$i]
    }
}

$newTimestamp = 44232.24014104684
$startRow = 380

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value2 = $names[$i]

    # Column B's visible text is the URL itself (matches every earlier row).
    $ws.Cells.Item($row, 2).Value2 = $cellTexts[$i]
    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $urls[$i], $subAddresses[$i]) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $urls[$i]) | Out-Null
    }
    # Hyperlinks.Add() stamps its own (slightly different) style xf, so
    # re-apply the exact "Hyperlink" cell style used by every other row
    # afterwards to keep the same style index (s="2") as the rest of the sheet.
    $ws.Cells.Item($row, 2).Style = "Hyperlink"

    $ws.Cells.Item($row, 3).Value2 = "Disponible"

    $ws.Cells.Item($row, 4).Value2 = $newTimestamp
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

"Updated rows 366-379 and appended rows 380-393"
